function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") '28.764.12'
$ws.Range("E2").Value = '  +3.09%  '

Set-TextValue $ws.Range("D3") '1.878.83'
$ws.Range("E3").Value = '  +3.14%  '

$ws.Range("E4").Value = '  +0.41%  '

Set-TextValue $ws.Range("D5") '324.59'
$ws.Range("E5").Value = '  -0.92%  '

Set-TextValue $ws.Range("D6") '1.004'
$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("E7").Value = '  +0.92%  '

Set-TextValue $ws.Range("D8") '0.3931'
$ws.Range("E8").Value = '  +2.29%  '

Set-TextValue $ws.Range("D9") '0.07928'
$ws.Range("E9").Value = '  +1.07%  '

Set-TextValue $ws.Range("D10") '0.9770'
$ws.Range("E10").Value = '  +1.90%  '

Set-TextValue $ws.Range("D11") '22.31'
$ws.Range("E11").Value = '  +2.13%  '

Set-TextValue $ws.Range("D12") '1.913.51'
$ws.Range("E12").Value = '  +4.39%  '

Set-TextValue $ws.Range("D13") '5.743'
$ws.Range("E13").Value = '  +1.78%  '

Set-TextValue $ws.Range("D14") '7.012'
$ws.Range("E14").Value = '  +2.39%  '

Set-TextValue $ws.Range("D15") '0.06960'
$ws.Range("E15").Value = '  +1.29%  '

Set-TextValue $ws.Range("D16") '88.65'
$ws.Range("E16").Value = '  +2.56%  '

$ws.Range("E17").Value = '  +0.45%  '

Set-TextValue $ws.Range("D18") '0.00001009'
$ws.Range("E18").Value = '  +1.60%  '

Set-TextValue $ws.Range("D19") '16.95'
$ws.Range("E19").Value = '  +1.89%  '

$ws.Range("E20").Value = '  +0.20%  '

Set-TextValue $ws.Range("D21") '28.776.67'
$ws.Range("E21").Value = '  +3.09%  '

Set-TextValue $ws.Range("D22") '5.363'
$ws.Range("E22").Value = '  +1.20%  '

Set-TextValue $ws.Range("D23") '11.09'
$ws.Range("E23").Value = '  +1.25%  '

Set-TextValue $ws.Range("D24") '2.118'
$ws.Range("E24").Value = '  +0.95%  '

Set-TextValue $ws.Range("D25") '2.115.44'
$ws.Range("E25").Value = '  +3.78%  '

Set-TextValue $ws.Range("D26") '153.57'
$ws.Range("E26").Value = '  +1.09%  '

Set-TextValue $ws.Range("D27") '19.39'
$ws.Range("E27").Value = '  +1.35%  '

Set-TextValue $ws.Range("D28") '5.749'
$ws.Range("E28").Value = '  -0.33%  '

Set-TextValue $ws.Range("D29") '2.000'
$ws.Range("E29").Value = '  +1.67%  '

Set-TextValue $ws.Range("D30") '119.98'
$ws.Range("E30").Value = '  +2.98%  '

Set-TextValue $ws.Range("D31") '0.09392'
$ws.Range("E31").Value = '  +1.82%  '

Set-TextValue $ws.Range("D32") '0.9393'
$ws.Range("E32").Value = '  +0.62%  '

Set-TextValue $ws.Range("D33") '5.309'
$ws.Range("E33").Value = '  +0.73%  '

Set-TextValue $ws.Range("D34") '1.354'
$ws.Range("E34").Value = '  +3.07%  '

Set-TextValue $ws.Range("D35") '3.348'
$ws.Range("E35").Value = '  +0.13%  '

Set-TextValue $ws.Range("D36") '0.05928'
$ws.Range("E36").Value = '  +0.02%  '

Set-TextValue $ws.Range("D37") '0.02117'
$ws.Range("E37").Value = '  -1.17%  '

Set-TextValue $ws.Range("D38") '1.158'
$ws.Range("E38").Value = '  +1.38%  '

Set-TextValue $ws.Range("D39") '7.900'
$ws.Range("E39").Value = '  +4.89%  '

Set-TextValue $ws.Range("D40") '0.5719'
$ws.Range("E40").Value = '  +2.72%  '

Set-TextValue $ws.Range("D41") '0.1797'
$ws.Range("E41").Value = '  +1.96%  '

Set-TextValue $ws.Range("D42") '9.981'
$ws.Range("E42").Value = '  +0.61%  '

Set-TextValue $ws.Range("D43") '0.07313'
$ws.Range("E43").Value = '  +4.67%  '

Set-TextValue $ws.Range("D44") '11.84'
$ws.Range("E44").Value = '  +2.50%  '

Set-TextValue $ws.Range("D45") '0.5340'
$ws.Range("E45").Value = '  +2.02%  '

Set-TextValue $ws.Range("D46") '1.148'
$ws.Range("E46").Value = '  -6.80%  '

Set-TextValue $ws.Range("D47") '2.113'
$ws.Range("E47").Value = '  -4.34%  '

Set-TextValue $ws.Range("D48") '1.844'
$ws.Range("E48").Value = '  +1.48%  '

Set-TextValue $ws.Range("D49") '114.12'
$ws.Range("E49").Value = '  +2.00%  '

Set-TextValue $ws.Range("D50") '2.376'
$ws.Range("E50").Value = '  +3.12%  '

Set-TextValue $ws.Range("D51") '1.004'
$ws.Range("E51").Value = '  +0.39%  '
